$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the row labels in column B (rows 4-29) -------------------------
# The regenerated notebook export spliced two new category names ("Holden"
# and "Rizzie Spiral") into the shared categorical list right after
# "Spiral5" / before "RotRing OmegaMax-90", and renamed "Thomas Hex" to
# "Matthies Hex". The row labels below reflect the label each row carries
# once the simulation was rerun with the updated category list.
$ws.Cells.Item(4,  2).Value = "Holden"
$ws.Cells.Item(5,  2).Value = "Rizzie Spiral"
$ws.Cells.Item(6,  2).Value = "RotRing OmegaMax-90"
$ws.Cells.Item(7,  2).Value = "Equal Angle"
$ws.Cells.Item(8,  2).Value = "Tilt Rotate"
$ws.Cells.Item(9,  2).Value = "CLR"
$ws.Cells.Item(10, 2).Value = "Rizzie Hex"
$ws.Cells.Item(11, 2).Value = "Matthies Hex"
$ws.Cells.Item(12, 2).Value = "Tilt Rotate_Partial"
$ws.Cells.Item(13, 2).Value = "RotRing OmegaMax-60"
$ws.Cells.Item(14, 2).Value = "Equal Angle_Partial"
$ws.Cells.Item(15, 2).Value = "Rizzie Hex_Partial"
$ws.Cells.Item(16, 2).Value = "ND Single"
$ws.Cells.Item(17, 2).Value = "RD Single"
$ws.Cells.Item(18, 2).Value = "TD Single"
$ws.Cells.Item(19, 2).Value = "Morris Single"
$ws.Cells.Item(20, 2).Value = "Ring Perpendicular to ND"
$ws.Cells.Item(21, 2).Value = "Ring Perpendicular to RD"
$ws.Cells.Item(22, 2).Value = "Ring Perpendicular to TD"
$ws.Cells.Item(23, 2).Value = "OffsetFTD"
$ws.Cells.Item(24, 2).Value = "OffsetATD"
$ws.Cells.Item(25, 2).Value = "OffsetF45"
$ws.Cells.Item(26, 2).Value = "OffsetA45"
$ws.Cells.Item(27, 2).Value = "OffsetFRD"
$ws.Cells.Item(28, 2).Value = "OffsetARD"
$ws.Cells.Item(29, 2).Value = "Gaussian Quadrature"

# --- Append two new simulation rows (30 and 31) -----------------------------
# Clone the formatting of the last existing data row (29) down into the two
# new rows so the bold/border style on column A survives, then overwrite the
# values for the new rows.
$ws.Range("A29:W29").Copy()
$ws.Range("A30:W30").PasteSpecial(-4122)
$ws.Range("A30:W30").PasteSpecial(-4123)

$ws.Range("A29:W29").Copy()
$ws.Range("A31:W31").PasteSpecial(-4122)
$ws.Range("A31:W31").PasteSpecial(-4123)

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "Michael-CCHex"
for ($col = 3; $col -le 23; $col++) {
    $ws.Cells.Item(30, $col).Value = 1
}

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "Michael-SNHex"
for ($col = 3; $col -le 23; $col++) {
    $ws.Cells.Item(31, $col).Value = 1
}
